$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price (D) and Volume(1h) (E) columns keep their original text formatting
# so numeric-looking strings (e.g. "1.001") are not coerced into Excel numbers/dates.
$ws.Range("D2:E51").NumberFormat = "@"

# Apply the updated cryptocurrency price / volume values
$ws.Range("D2").Value = "29.811.80"
$ws.Range("E2").Value = "  +8.23%  "
$ws.Range("D3").Value = "1.955.11"
$ws.Range("E3").Value = "  +6.58%  "
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  -0.26%  "
$ws.Range("D5").Value = "342.05"
$ws.Range("E5").Value = "  +2.58%  "
$ws.Range("E6").Value = "  -0.22%  "
$ws.Range("D7").Value = "0.4779"
$ws.Range("E7").Value = "  +3.95%  "
$ws.Range("D8").Value = "0.4146"
$ws.Range("E8").Value = "  +8.25%  "
$ws.Range("D9").Value = "47.96"
$ws.Range("E9").Value = "  +3.54%  "
$ws.Range("D10").Value = "0.08250"
$ws.Range("E10").Value = "  +5.02%  "
$ws.Range("D11").Value = "1.038"
$ws.Range("E11").Value = "  +7.81%  "
$ws.Range("D12").Value = "22.81"
$ws.Range("E12").Value = "  +8.09%  "
$ws.Range("D13").Value = "1.946.65"
$ws.Range("E13").Value = "  +7.16%  "
$ws.Range("D14").Value = "6.178"
$ws.Range("E14").Value = "  +5.56%  "
$ws.Range("D15").Value = "7.407"
$ws.Range("E15").Value = "  +4.67%  "
$ws.Range("D16").Value = "92.09"
$ws.Range("E16").Value = "  +2.62%  "
$ws.Range("D17").Value = "1.002"
$ws.Range("E17").Value = "  -0.21%  "
$ws.Range("E18").Value = "  +3.54%  "
$ws.Range("D19").Value = "0.06696"
$ws.Range("E19").Value = "  +1.89%  "
$ws.Range("E20").Value = "  +5.07%  "
$ws.Range("D21").Value = "1.0000"
$ws.Range("E21").Value = "  -0.39%  "
$ws.Range("D22").Value = "29.778.90"
$ws.Range("E22").Value = "  +8.24%  "
$ws.Range("E23").Value = "  +5.35%  "
$ws.Range("E24").Value = "  +4.42%  "
$ws.Range("E25").Value = "  +0.39%  "
$ws.Range("D26").Value = "2.177.28"
$ws.Range("E26").Value = "  +6.85%  "
$ws.Range("D27").Value = "161.34"
$ws.Range("E27").Value = "  +1.67%  "
$ws.Range("D28").Value = "20.25"
$ws.Range("E28").Value = "  +4.38%  "
$ws.Range("D29").Value = "2.183"
$ws.Range("E29").Value = "  +6.71%  "
$ws.Range("D30").Value = "5.697"
$ws.Range("E30").Value = "  +7.72%  "
$ws.Range("D31").Value = "123.19"
$ws.Range("E31").Value = "  +4.23%  "
$ws.Range("E32").Value = "  +7.97%  "
$ws.Range("D33").Value = "0.09642"
$ws.Range("E33").Value = "  +2.58%  "
$ws.Range("D34").Value = "1.479"
$ws.Range("E34").Value = "  +12.15%  "
$ws.Range("E35").Value = "  +3.17%  "
$ws.Range("D36").Value = "5.521"
$ws.Range("E36").Value = "  +5.72%  "
$ws.Range("D37").Value = "0.06281"
$ws.Range("E37").Value = "  +5.57%  "
$ws.Range("D38").Value = "0.02316"
$ws.Range("E38").Value = "  +5.96%  "
$ws.Range("D39").Value = "8.478"
$ws.Range("E39").Value = "  +4.23%  "
$ws.Range("D40").Value = "0.6095"
$ws.Range("E40").Value = "  +6.19%  "
$ws.Range("D41").Value = "1.185"
$ws.Range("E41").Value = "  +3.58%  "
$ws.Range("D42").Value = "10.75"
$ws.Range("E42").Value = "  +7.87%  "
$ws.Range("E43").Value = "  +3.94%  "
$ws.Range("E44").Value = "  -0.29%  "
$ws.Range("D45").Value = "1.280"
$ws.Range("E45").Value = "  +0.96%  "
$ws.Range("D46").Value = "2.402"
$ws.Range("E46").Value = "  +35.16%  "
$ws.Range("D47").Value = "0.5719"
$ws.Range("E47").Value = "  +6.12%  "
$ws.Range("D48").Value = "12.51"
$ws.Range("E48").Value = "  +6.30%  "
$ws.Range("D49").Value = "0.07385"
$ws.Range("E49").Value = "  +7.95%  "
$ws.Range("D50").Value = "1.989"
$ws.Range("E50").Value = "  +4.43%  "
$ws.Range("D51").Value = "113.34"
$ws.Range("E51").Value = "  +1.97%  "

# Row 45/46 swapped coins (RenderToken <-> WEMIXTOKEN) together with their data
$ws.Range("B45").Value = "WEMIXTOKEN"
$ws.Range("C45").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("B46").Value = "RenderToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"

# Row 49/50 swapped coins (NEARProtocol <-> Cronos) together with their data
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("B50").Value = "NEARProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
